# Weekly update: a new daily quotation row is inserted at the top of the
# "Papa" (potato) price history for Agricola del Norte S.A. de Arica,
# pushing every existing record down by one row (row 23 -> 24, 24 -> 25, ...,
# 71 -> 72). The sheet's used range grows from A1:R71 to A1:R72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 23 (the oldest-until-now
# entry). Excel shifts rows 23..71 down to 24..72 and extends the used
# range accordingly.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the latest market quotation.
$ws.Range("A23").Value = 1
$ws.Range("B23").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C23").Value = "Arica y Parinacota"
$ws.Range("D23").Value2 = 44707
$ws.Range("E23").Value = 15
$ws.Range("F23").Value = 100114001
$ws.Range("G23").Value = "Papa"
$ws.Range("H23").Value = "Asterix"
$ws.Range("I23").Value = "1a (cosecha lavada)"
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 9500
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = 9750
$ws.Range("N23").Value = "`$/malla 25 kilos"
$ws.Range("O23").Value = "Región de Los Lagos"
$ws.Range("P23").Value = 390
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
